$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B for TestDescription, shifting old Execute/InvocationCount columns right
$ws.Range("B1").EntireColumn.Insert()

# Header row
$ws.Range("A1").Value = "TestName"
$ws.Range("B1").Value = "TestDescription"
$ws.Range("C1").Value = "Execute"
$ws.Range("D1").Value = "InvocationCount"

# Row 2
$ws.Range("A2").Value = "login"
$ws.Range("B2").Value = "This is a test to validate login in to OrangeHrm"
$ws.Range("C2").Value = "Yes"
$ws.Range("D2").Value = "'2"

# Row 3
$ws.Range("A3").Value = "login"
$ws.Range("B3").Value = "This is a copy "
$ws.Range("C3").Value = "No"
$ws.Range("D3").Value = "'1"

# Column widths (leave column A untouched to preserve its bestFit width)
$ws.Columns.Item(2).ColumnWidth = 18.85

# Selection
$ws.Range("A6").Select()
